# Weekly fruit/vegetable price update:
# Insert a new record (row) before the current row 95, shifting the
# existing rows 95-103 down to 96-104, and populate the new row with the
# latest "Locoto" price observation for "Vega Modelo de Temuco".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 95:103 down by one row, opening up a blank row 95.
$ws.Rows("95:95").Insert()

# Fill in the new row 95 with the new weekly observation.
$ws.Range("A95").Value = 10
$ws.Range("B95").Value = "Vega Modelo de Temuco"
$ws.Range("C95").Value = "La Araucanía"
$ws.Range("D95").Value = 45194
$ws.Range("E95").Value = 9
$ws.Range("F95").Value = 100112042
$ws.Range("G95").Value = "Locoto"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 100
$ws.Range("K95").Value = 2200
$ws.Range("L95").Value = 2200
$ws.Range("M95").Value = 2200
$ws.Range("N95").Value = "$/kilo"
$ws.Range("O95").Value = "Región de Arica y Parinacota"
$ws.Range("P95").Value = 2200
$ws.Range("Q95").Value = 1
$ws.Range("R95").Value = "Hortaliza"
